# B1--and-B2-PowerPoint.pptx edit
# 1) Slide 5's table switches to a different built-in table style.
# 2) The presentation's theme colour scheme is swapped from the
#    "Integral" (Red Violet) palette to the "Office" (Office Theme)
#    palette.

$p = $ppt.ActivePresentation

# --- 1. Update the table style on slide 5 -----------------------------
$slide = $p.Slides.Item(5)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$tableShape.Table.ApplyStyle("{D7C19E04-D30B-43C5-8AB9-A267CC960171}")

# --- 2. Swap the theme colour scheme (Integral/Red Violet -> Office) --
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,         # dk1    000000
    16777215,  # lt1    FFFFFF
    6968388,   # dk2    44546A
    15132391,  # lt2    E7E6E6
    13998939,  # accent1 5B9BD5
    3243501,   # accent2 ED7D31
    10855845,  # accent3 A5A5A5
    49407,     # accent4 FFC000
    12874308,  # accent5 4472C4
    4697456,   # accent6 70AD47
    12673797,  # hlink   0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
